$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lijst")
$st = $wb.Styles.Add("MyBorderStyle")
$st.Borders(7).LineStyle = 1
$st.Borders(10).LineStyle = 1
$rng = $ws.Range("J2:J126")
$rng.Value = "ja"
$rng.Style = "MyBorderStyle"
Write-Host "ok"
